$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 20:33"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1200874
$ws.Range("C4").Value = 12752
$ws.Range("D4").Value = 181602
$ws.Range("E4").Value = 950151
$ws.Range("G4").Value = 523
$ws.Range("H4").Value = 69121

# Row 8 - Francia
$ws.Range("B8").Value = 169462
$ws.Range("C8").Value = 769
$ws.Range("D8").Value = 51371
$ws.Range("E8").Value = 92890
$ws.Range("F8").Value = 3696
$ws.Range("G8").Value = 306
$ws.Range("H8").Value = 25201

# Row 62 - Barein
$ws.Range("B62").Value = 3533
$ws.Range("C62").Value = 150
$ws.Range("E62").Value = 1807

# Row 169 - Libia
$ws.Range("D169").Value = 23
$ws.Range("E169").Value = 37

# Row 178 - San Martin (Parte Francesa)
$ws.Range("D178").Value = 29
$ws.Range("E178").Value = 6
$ws.Range("F178").Value = 1
